$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "custom accuracy": round the numeric sample values in row 5 (columns B:AH)
# down to 2 decimal places.
$rng = $ws.Range("B5:AH5")
foreach ($cell in $rng.Cells) {
    $cell.Value2 = [Math]::Round($cell.Value2, 2)
}

# A couple of values sit exactly on a rounding midpoint and the source data
# set was produced with "round half to even" (banker's rounding) rather than
# "round half away from zero" - correct those two cells explicitly so the
# stored values match the originally published data set exactly.
$ws.Range("G5").Value2 = 0.62
$ws.Range("P5").Value2 = 1.21

# "데이터 1000개" - the sample data set was trimmed; drop the last sample row.
$ws.Rows.Item(6).Delete()
